# Re-run the Area3 variogram-parameter optimisation with a smaller
# search radius (110000 -> 75000) and paste in the new model output,
# then tidy up the "Area3" label and highlight its header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the Area3 label (drop the trailing space) -----------------
$ws.Range("U2").Value = "Area3"

# --- 2. New search_radius input for the Area3 table (V3:V11) ----------
$ws.Range("V3:V11").Value = 75000

# --- 3. New optimisation output pasted into X, Y, Z, AB (rows 3-11) ---
$xVals = @(0.38033031902469272, 0.69008680305220649, 0.74321284253246644, 0.80147713775890628, 0.84652785808803543, 0.88163691075563233, 0.91648760976161969, 0.9304421052752393, 0.93761111150821319)
$yVals = @(3.6025795623284291, 3.3636685309936909, 3.3523063974154081, 3.3244513266097129, 3.3135987201128718, 3.3051287071358479, 3.2800970521033079, 3.2725045109609678, 3.268319359703324)
$zVals = @(0.75180000000000002, 1.3849, 1.7165999999999999, 1.95, 2.0798000000000001, 2.2050000000000001, 2.3134000000000001, 2.3569, 2.37229)

for ($i = 0; $i -lt 9; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 24).Value = $xVals[$i]   # column X
    $ws.Cells.Item($row, 25).Value = $yVals[$i]   # column Y
    $ws.Cells.Item($row, 26).Value = $zVals[$i]   # column Z
}

# AB3:AB11 is a constant repeated down the column
$ws.Range("AB3:AB11").Value = 3.748694078222786

# --- 4. Bold/highlight the Area3 table header row ----------------------
$ws.Range("U2:AB2").Font.Bold = $true
$ws.Range("U3:U11").Font.Bold = $true

# --- 5. Update the saved view state (scroll position + selection) ------
$ws.Range("AB23").Select()

$wb.Save()
